# Add a new "BB" column to the right of the existing "BA" column.
#
# Row 1 (the date header row) gets a new date, copied-formatted from BA1.
# Rows 3-18 repeat (carry forward) the same value that is already in the
# BA column of that row.
# Rows 19-21 get a new, distinct forecast value.
# Rows 2 and 22 are untouched (they only have a value in column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: new date header, keep the same (bold/centered/date) formatting as BA1.
$ws.Range("BA1").Copy($ws.Range("BB1"))
$ws.Range("BB1").Value = 45986

# Rows 3-18: carry the BA value forward into BB (same formatting - default).
for ($r = 3; $r -le 18; $r++) {
    $ws.Range("BA$r").Copy($ws.Range("BB$r"))
}

# Rows 19-21: new forecast values, distinct from the BA column.
$ws.Range("BB19").Value = 0.8976398032236155
$ws.Range("BB20").Value = 0.7456737245741252
$ws.Range("BB21").Value = 0.7805163230192314
